# End of the dispersal draft / start of the auction:
#  - fix a few rows whose "position" value was recorded incorrectly
#  - append the newly drafted/auctioned players

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("draftpicks")

# --- Correct previously mis-recorded positions -----------------------
$ws.Cells.Item(35, 4).Value  = "SS"
$ws.Cells.Item(36, 4).Value  = "3B"
$ws.Cells.Item(125, 4).Value = "OF"
$ws.Cells.Item(162, 4).Value = "SS"

# --- Append new rows for the latest picks -----------------------------
# (team, player, salary, position, draftedSerial) for rows 199..204
$newRows = @{
    199 = @("dsb",     "Lourdes Gurriel Jr.", 7,  "2B", 43465)
    200 = @("balco",   "Ross Stripling",      5,  "P",  43103)
    201 = @("chicago", "Jorge Polanco",       13, "SS", 43103)
    202 = @("deano",   "Michael Kopech",      0,  "B",  43103)
    203 = @("dsb",     "Niko Goodrum",        5,  "2B", 43103)
    204 = @("balco",   "Jose Alvarado",       5,  "P",  43103)
}

# Populate the player-name (column B) cells in the same order the picks
# were made so new shared-string entries land in the right sequence.
$order = @(200, 201, 202, 203, 204, 199)
foreach ($r in $order) {
    $ws.Cells.Item($r, 2).Value = $newRows[$r][1]
}

for ($r = 199; $r -le 204; $r++) {
    $data = $newRows[$r]

    $ws.Cells.Item($r, 1).Value = $data[0]

    # carry over the currency/date number formatting from the row above
    # without introducing new style/numFmt entries
    $ws.Cells.Item($r - 1, 3).Copy()
    $ws.Cells.Item($r, 3).PasteSpecial(-4122)
    $ws.Cells.Item($r - 1, 5).Copy()
    $ws.Cells.Item($r, 5).PasteSpecial(-4122)
    $excel.CutCopyMode = $false

    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
}

# --- Reflect the newly active selection --------------------------------
$ws.Range("B199").Select()
try { $excel.ActiveWindow.ScrollRow = 185 } catch {}
